$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update FlowScriptFile (K2) and Time (I2) values
$ws.Range("K2").Value = "ScriptedFlow_POManagementScript_PrePROD_50bde6f4.xlsx"
$ws.Range("I2").Value = "'20"

# Update the view: scroll back to top-left and move the selection to G9
$ws.Range("A1").Select()
$ws.Range("G9").Select()
